# Delete row 641 ("「珈琲で一番大事な事！誰と飲むか？」" post) and shift
# subsequent rows up, matching the author's commit that removed this post.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(641).Delete()
